$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header fields: query -> prompt, positive -> chosen, negative -> rejected
$ws.Range("A1").Value = "prompt"
$ws.Range("B1").Value = "chosen"
$ws.Range("C1").Value = "rejected"

# Apply a white fill across the whole used data range (A1:E10), matching the
# fillId/applyFill additions seen on the data-row cell styles.
$ws.Range("A1:E10").Interior.Color = 16777215
